# new update with respect to AAP and LAF module
# Insert a "SortByoption" worksheet between "AKmodule" and "Kaif Khan",
# populate it with the sort-option test data, and refresh the
# AKmodule sheet's saved selection.

$wb = $excel.ActiveWorkbook

$akModule = $wb.Worksheets.Item(1)

# Insert the new sheet right after AKmodule (i.e. before the 2nd sheet,
# "Kaif Khan"), matching sheets order AKmodule / SortByoption / Kaif Khan.
$sortSheet = $wb.Worksheets.Add($null, $akModule)
$sortSheet.Name = "SortByoption"

# --- Values -------------------------------------------------------------
$sortSheet.Range("A1").Value = "TC_ID"
$sortSheet.Range("B1").Value = "SortOption"

$sortSheet.Range("A2").Value = 1
$sortSheet.Range("B2").Value = " Low to High price"

$sortSheet.Range("A3").Value = 2
$sortSheet.Range("B3").Value = " Highest to Lowest price"

$sortSheet.Range("A4").Value = 3
$sortSheet.Range("B4").Value = "Low to High discount"

$sortSheet.Range("A5").Value = 4
$sortSheet.Range("B5").Value = "Highest to Lowest discount"

# --- Formatting -----------------------------------------------------------
# Header row: bold, centered both ways, wrap text. Build the format once on
# A1 and fan it out via format-only paste so we don't leave behind
# intermediate/unused style entries.
$sortSheet.Range("A1").Font.Bold = $true
$sortSheet.Range("A1").VerticalAlignment = -4108
$sortSheet.Range("A1").HorizontalAlignment = -4108
$sortSheet.Range("A1").WrapText = $true
$sortSheet.Range("A1").Copy()
$sortSheet.Range("B1").PasteSpecial(-4122)

# Data rows: vertical-centered, wrap text (default font, no horizontal
# centering). Build on A2 and fan out the same way.
$sortSheet.Range("A2").VerticalAlignment = -4108
$sortSheet.Range("A2").WrapText = $true
$sortSheet.Range("A2").Copy()
$sortSheet.Range("B2:B5").PasteSpecial(-4122)
$sortSheet.Range("A3:A5").PasteSpecial(-4122)

# --- Row heights ------------------------------------------------------
$sortSheet.Rows.Item(1).RowHeight = 30
$sortSheet.Rows.Item(2).RowHeight = 45
$sortSheet.Rows.Item(3).RowHeight = 60
$sortSheet.Rows.Item(4).RowHeight = 45
$sortSheet.Rows.Item(5).RowHeight = 60

# --- Selections / active sheet -----------------------------------------
# AKmodule keeps a plain selection at B56 once it stops being the active tab.
[void]$akModule.Range("B56").Select()

# SortByoption becomes the active/selected sheet.
[void]$sortSheet.Activate()
[void]$sortSheet.Range("A1:B5").Select()
